# "derniers retour powerpoint": the récapitulatif table on slide 14 showed
# the vacance rate placeholders with a literal "%" baked into the template
# text (the value supplied at merge time already carries the percent sign),
# so the trailing "%" is stripped from each of the three scenario cells.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$tbl = $s.Shapes.Item(2).Table

# Row 5 ("Logements vacants > 2 ans"), columns 2-4 hold {{vacance1}},
# {{vacance2}} and {{vacance3}} for the three scenarios.
$tbl.Cell(5, 2).Shape.TextFrame.TextRange.Text = "{{vacance1}} "
$tbl.Cell(5, 3).Shape.TextFrame.TextRange.Text = "{{vacance2}} "
$tbl.Cell(5, 4).Shape.TextFrame.TextRange.Text = "{{vacance3}}"
